# Update "想去人数" (want-to-go count) figures across the four sheets
# to reflect the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6
$ws1.Range("F5").Value = 19489
$ws1.Range("F7").Value = 2205
$ws1.Range("F10").Value = 427
$ws1.Range("F11").Value = 685
$ws1.Range("F17").Value = 261

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 284
$ws2.Range("F8").Value = 128
$ws2.Range("F12").Value = 9

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 639

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 639
$ws4.Range("F6").Value = 6
$ws4.Range("F10").Value = 19489
$ws4.Range("F15").Value = 284
$ws4.Range("F16").Value = 2205
$ws4.Range("F18").Value = 128
$ws4.Range("F20").Value = 427
$ws4.Range("F21").Value = 685
$ws4.Range("F30").Value = 9
$ws4.Range("F31").Value = 261
